# Apply updated cryptocurrency price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.276.84'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").Value = '3.070.16'
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'574.65"
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("D6").Value = "'170.07"
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '3.067.12'
$ws.Range("E8").Value = '  -1.69%  '
$ws.Range("E9").Value = '  -2.30%  '
$ws.Range("D10").Value = "'6.26"
$ws.Range("E10").Value = '  -2.02%  '
$ws.Range("E11").Value = '  -2.75%  '
$ws.Range("D12").Value = "'0.467"
$ws.Range("E12").Value = '  -3.15%  '
$ws.Range("E13").Value = '  -3.98%  '
$ws.Range("D14").Value = "'35.62"
$ws.Range("E14").Value = '  -4.24%  '
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = '3.582.18'
$ws.Range("D17").Value = '66.235.45'
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").Value = "'6.93"
$ws.Range("E18").Value = '  -3.13%  '
$ws.Range("D19").Value = '3.071.38'
$ws.Range("E19").Value = '  -1.58%  '
$ws.Range("D20").Value = "'16.57"
$ws.Range("E20").Value = '  +1.66%  '
$ws.Range("E21").Value = '  +1.68%  '
$ws.Range("D22").Value = "'0.683"
$ws.Range("E22").Value = '  -3.68%  '
$ws.Range("D23").Value = "'7.64"
$ws.Range("E23").Value = '  -2.86%  '
$ws.Range("D24").Value = "'82.20"
$ws.Range("E24").Value = '  -2.03%  '
$ws.Range("E25").Value = '  -4.55%  '
$ws.Range("E26").Value = '  -3.65%  '
$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("E30").Value = '  -5.51%  '
$ws.Range("D31").Value = "'2.58"
$ws.Range("E31").Value = '  -3.64%  '
$ws.Range("D32").Value = "'27.61"
$ws.Range("E32").Value = '  -3.53%  '
$ws.Range("E33").Value = '  -3.32%  '
$ws.Range("D34").Value = '0.0₃0916'
$ws.Range("E34").Value = '  -3.70%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = "'47.94"
$ws.Range("E36").Value = '  +1.99%  '
$ws.Range("E37").Value = '  -3.33%  '
$ws.Range("E38").Value = '  -5.05%  '
$ws.Range("E39").Value = '  -1.20%  '
$ws.Range("E41").Value = '  -5.15%  '
$ws.Range("E42").Value = '  -4.53%  '
$ws.Range("D43").Value = '2.775.40'
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("E45").Value = '  -2.82%  '
$ws.Range("D46").Value = "'134.43"
$ws.Range("E46").Value = '  -1.14%  '
$ws.Range("D47").Value = "'364.25"
$ws.Range("E47").Value = '  -4.88%  '
$ws.Range("D49").Value = "'24.16"
$ws.Range("E49").Value = '  -3.34%  '
$ws.Range("E50").Value = '  -2.71%  '
$ws.Range("D51").Value = "'0.105"
$ws.Range("E51").Value = '  -2.42%  '
